$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.048.84'
$ws.Range("E2").Value = '  -4.61%  '
$ws.Range("D3").Value = '3.178.02'
$ws.Range("E3").Value = '  -5.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '532.06'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -6.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.37'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -8.19%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '3.176.07'
$ws.Range("E8").Value = '  -5.25%  '
$ws.Range("E9").Value = '  -6.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.25'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -8.42%  '
$ws.Range("E11").Value = '  -8.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.394'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -4.87%  '
$ws.Range("D13").Value = '3.728.42'
$ws.Range("E13").Value = '  -4.94%  '
$ws.Range("E14").Value = '  -0.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.66'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -7.18%  '
$ws.Range("D16").Value = '3.184.80'
$ws.Range("E16").Value = '  -4.99%  '
$ws.Range("D17").Value = '58.178.66'
$ws.Range("E17").Value = '  -4.40%  '
$ws.Range("E18").Value = '  -8.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.84'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -7.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.22'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -8.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.07'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -9.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '358.36'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.70'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -6.81%  '
$ws.Range("E25").Value = '  -7.90%  '
$ws.Range("D26").Value = '3.320.12'
$ws.Range("E26").Value = '  -5.11%  '
$ws.Range("E27").Value = '  -3.89%  '
$ws.Range("E28").Value = '  -12.14%  '
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.90'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.94%  '
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("E32").Value = '  -8.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.95'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -9.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '21.63'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -5.23%  '
$ws.Range("E35").Value = '  -7.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.93'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -7.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.40'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.93%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.43'
$ws.Range("D38").ClearFormats()
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.27'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -8.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.86'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -9.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0702'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -6.50%  '
$ws.Range("D42").Value = '3.210.69'
$ws.Range("E42").Value = '  -5.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.55'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.30%  '
$ws.Range("E44").Value = '  -7.01%  '
$ws.Range("E45").Value = '  -6.83%  '
$ws.Range("E46").Value = '  -4.39%  '
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("E48").Value = '  -8.32%  '
$ws.Range("D49").Value = '2.275.18'
$ws.Range("E49").Value = '  -7.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.19'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -6.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.54'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -8.08%  '
